$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "harvester" column (B) from "Retrofitted_0731" to "S.GISH" for data rows 2-19
# and populate the new "experimentDesign" column (D) with "90minuteInduction" for the same rows.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Reflect the selection change recorded in the sheet view (active cell D3, selection D3:D19)
$ws.Range("D3:D19").Select()
